$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-21 Thursday", "2024-11-22 Friday"),
    @("530×8=", "921×4="),
    @("743×9=", "289×3="),
    @("797×2=", "148×9="),
    @("877×8=", "152×2="),
    @("209×3=", "956×2="),
    @("229×9=", "691×7="),
    @("382×9=", "601×4="),
    @("902×6=", "267×4="),
    @("339×3=", "665×3="),
    @("597×5=", "677×9="),
    @("145×2=", "647×7="),
    @("528×8=", "990×4="),
    @("139×2=", "221×9="),
    @("549×6=", "539×2="),
    @("161×7=", "709×8="),
    @("551×3=", "383×9="),
    @("140×3=", "638×3="),
    @("657×8=", "754×6="),
    @("811×4=", "132×8="),
    @("163×5=", "922×5="),
    @("855×6=", "221×2="),
    @("154×3=", "741×5="),
    @("881×6=", "137×3="),
    @("181×9=", "355×6="),
    @("309×2=", "811×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
